# Updated symbol list on Sun Dec 18 09:21:44 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Price cells (column D) hold numeric-looking text (e.g. trailing zeros like
# '1.070' or '0.3290') that must stay text, so mark them Text before writing.
$priceCells = @(
    'D2', 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15',
    'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27',
    'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D47', 'D48', 'D49', 'D50'
)
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '247.55'
$ws.Range('D4').Value = '5.565'
$ws.Range('D5').Value = '0.05624'
$ws.Range('D6').Value = '3.395'
$ws.Range('D7').Value = '6.475'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = '0.8033'
$ws.Range('E8').Value = '7MXTokenMX'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').Value = '1.070'
$ws.Range('E9').Value = '8FTXTokenFTT'
$ws.Range('D10').Value = '0.1432'
$ws.Range('D11').Value = '0.07320'
$ws.Range('D12').Value = '0.03203'
$ws.Range('D13').Value = '0.02992'
$ws.Range('D14').Value = '0.09256'
$ws.Range('D15').Value = '0.001676'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').Value = '2.968'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = '0.04695'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').Value = '0.0006039'
$ws.Range('E18').Value = '17OneONE'
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D19').Value = '0.006275'
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('B20').Value = 'BitKan'
$ws.Range('C20').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D20').Value = '0.001052'
$ws.Range('E20').Value = '19BitKanKAN'
$ws.Range('B21').Value = 'HotbitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D21').Value = '0.003828'
$ws.Range('E21').Value = '20HotbitTokenHTB'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D22').Value = '0.0001502'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('B23').Value = 'UpBots'
$ws.Range('C23').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D23').Value = '0.0004006'
$ws.Range('E23').Value = '22UpBotsUBXT'
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').Value = '3.986'
$ws.Range('E24').Value = '23LEOLEO'
$ws.Range('B25').Value = 'BTSEToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D25').Value = '2.113'
$ws.Range('E25').Value = '24BTSETokenBTSE'
$ws.Range('B26').Value = 'BitpandaEcosystemToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D26').Value = '0.3290'
$ws.Range('E26').Value = '25BitpandaEcosystemTokenBEST'
$ws.Range('B27').Value = 'ProBitToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D27').Value = '0.1292'
$ws.Range('E27').Value = '26ProBitTokenPROBBestin24h'
$ws.Range('D40').Value = '0.04186'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').Value = '0.1046'
$ws.Range('E41').Value = '40BKEXTokenBKK'
$ws.Range('D42').Value = '0.002975'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').Value = '0.003236'
$ws.Range('E43').Value = '42KickTokenKICKWorstin24h'
$ws.Range('D44').Value = '0.009619'
$ws.Range('D45').Value = '0.00005654'
$ws.Range('D47').Value = '0.6810'
$ws.Range('D48').Value = '0.02688'
$ws.Range('E48').Value = '47BOLOBOLO'
$ws.Range('D49').Value = '0.00002103'
$ws.Range('D50').Value = '0.01012'
